$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.601557195186615
$ws.Range("B1").Value = 0.5433940887451172
$ws.Range("C1").Value = 3.65119457244873
$ws.Range("D1").Value = 1.706456303596497
$ws.Range("E1").Value = 1.191561937332153
